$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.058283626036012
$ws.Range("D2").Value = 1.057068940049686
$ws.Range("E2").Value = 1.071917139059091
$ws.Range("F2").Value = 1.079183756150517
$ws.Range("I2").Value = 1.053044569886574
$ws.Range("J2").Value = 1.063275649752193
$ws.Range("K2").Value = 1.059804921180433
$ws.Range("L2").Value = 1.074613038059806
$ws.Range("M2").Value = 1.081860473541999
$ws.Range("N2").Value = 1.024713022826894

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.059432031319226
$ws.Range("D3").Value = 1.057944592456873
$ws.Range("E3").Value = 1.073056454274037
$ws.Range("F3").Value = 1.080397916716307
$ws.Range("I3").Value = 1.053440013948049
$ws.Range("J3").Value = 1.064075903650775
$ws.Range("K3").Value = 1.060494179694598
$ws.Range("L3").Value = 1.075568202058651
$ws.Range("M3").Value = 1.082891692691552
$ws.Range("N3").Value = 1.024991502360548

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.06017494666825
$ws.Range("D4").Value = 1.058510966221081
$ws.Range("E4").Value = 1.073793849110387
$ws.Range("F4").Value = 1.08118387992226
$ws.Range("I4").Value = 1.053694514040123
$ws.Range("J4").Value = 1.064592971287492
$ws.Range("K4").Value = 1.06093930496347
$ws.Range("L4").Value = 1.076185849951639
$ws.Range("M4").Value = 1.083558704693529
$ws.Range("N4").Value = 1.02517117764816

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.060487226575748
$ws.Range("D5").Value = 1.058749014571406
$ws.Range("E5").Value = 1.07410389428719
$ws.Range("F5").Value = 1.081514376526685
$ws.Range("I5").Value = 1.053801176207149
$ws.Range("J5").Value = 1.064810167627041
$ws.Range("K5").Value = 1.061126227583128
$ws.Range("L5").Value = 1.076445412477816
$ws.Range("M5").Value = 1.083839056210711
$ws.Range("N5").Value = 1.02524658885322

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.060539657338726
$ws.Range("D6").Value = 1.058788980710869
$ws.Range("E6").Value = 1.074155954865353
$ws.Range("F6").Value = 1.081569872945652
$ws.Range("I6").Value = 1.05381906593335
$ws.Range("J6").Value = 1.064846625377814
$ws.Range("K6").Value = 1.0611576005281
$ws.Range("L6").Value = 1.076488988524149
$ws.Range("M6").Value = 1.083886124952481
$ws.Range("N6").Value = 1.025259243441272

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.060179119529886
$ws.Range("D7").Value = 1.058514147251188
$ws.Range("E7").Value = 1.073797991773685
$ws.Range("F7").Value = 1.081188295727562
$ws.Range("I7").Value = 1.053695940560123
$ws.Range("J7").Value = 1.064595874177678
$ws.Range("K7").Value = 1.060941803449511
$ws.Range("L7").Value = 1.076189318617561
$ws.Range("M7").Value = 1.083562451000144
$ws.Range("N7").Value = 1.025172185784595

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.058671772963881
$ws.Range("D8").Value = 1.057364919099344
$ws.Range("E8").Value = 1.072302138879742
$ws.Range("F8").Value = 1.079594021734451
$ws.Range("I8").Value = 1.053178497748706
$ws.Range("J8").Value = 1.063546255043376
$ws.Range("K8").Value = 1.06003803998313
$ws.Range("L8").Value = 1.074935924981564
$ws.Range("M8").Value = 1.082209032374974
$ws.Range("N8").Value = 1.024807243937719

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.056014201032121
$ws.Range("D9").Value = 1.055338040625823
$ws.Range("E9").Value = 1.06966760656963
$ws.Range("F9").Value = 1.076787112510216
$ws.Range("I9").Value = 1.052256121220701
$ws.Range("J9").Value = 1.061690915490965
$ws.Range("K9").Value = 1.058438799351482
$ws.Range("L9").Value = 1.072724130918441
$ws.Range("M9").Value = 1.079822136643628
$ws.Range("N9").Value = 1.024160185417292

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.054241437317015
$ws.Range("D10").Value = 1.053985559207882
$ws.Range("E10").Value = 1.067912101718148
$ws.Range("F10").Value = 1.074917389003522
$ws.Range("I10").Value = 1.051634067082671
$ws.Range("J10").Value = 1.060450092673482
$ws.Range("K10").Value = 1.05736810471349
$ws.Range("L10").Value = 1.071247427199458
$ws.Range("M10").Value = 1.078229469531714
$ws.Range("N10").Value = 1.023726126345329

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.053473542661607
$ws.Range("D11").Value = 1.053399621989384
$ws.Range("E11").Value = 1.067152136429127
$ws.Range("F11").Value = 1.074108129605765
$ws.Range("I11").Value = 1.051363011933788
$ws.Range("J11").Value = 1.059911859670331
$ws.Range("K11").Value = 1.056903398505288
$ws.Range("L11").Value = 1.070607469772518
$ws.Range("M11").Value = 1.077539480288031
$ws.Range("N11").Value = 1.023537534644838

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.053188269104125
$ws.Range("D12").Value = 1.053181932524826
$ws.Range("E12").Value = 1.066869877208678
$ws.Range("F12").Value = 1.07380758472052
$ws.Range("I12").Value = 1.051262073733321
$ws.Range("J12").Value = 1.059711792425879
$ws.Range("K12").Value = 1.056730621621374
$ws.Range("L12").Value = 1.070369679524328
$ws.Range("M12").Value = 1.077283133091768
$ws.Range("N12").Value = 1.023467386726176

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.053249463238329
$ws.Range("D13").Value = 1.053228629759162
$ws.Range("E13").Value = 1.066930421626419
$ws.Range("F13").Value = 1.073872050375814
$ws.Range("I13").Value = 1.051283736925471
$ws.Range("J13").Value = 1.05975471404119
$ws.Range("K13").Value = 1.056767690297467
$ws.Range("L13").Value = 1.070420690049944
$ws.Range("M13").Value = 1.077338122921364
$ws.Range("N13").Value = 1.023482438066153

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.053449962753826
$ws.Range("D14").Value = 1.053381628658428
$ws.Range("E14").Value = 1.067128804276658
$ws.Range("F14").Value = 1.074083285458576
$ws.Range("I14").Value = 1.051354673581546
$ws.Range("J14").Value = 1.059895324983703
$ws.Range("K14").Value = 1.056889120070569
$ws.Range("L14").Value = 1.070587815627237
$ws.Range("M14").Value = 1.077518291667522
$ws.Range("N14").Value = 1.023531738167162

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.053573491332615
$ws.Range("D15").Value = 1.05347589017631
$ws.Range("E15").Value = 1.067251037744813
$ws.Range("F15").Value = 1.074213440957553
$ws.Range("I15").Value = 1.051398346022993
$ws.Range("J15").Value = 1.059981940974401
$ws.Range("K15").Value = 1.056963915165446
$ws.Range("L15").Value = 1.070690776367321
$ws.Range("M15").Value = 1.077629292321613
$ws.Range("N15").Value = 1.023562100784175

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.054292393983231
$ws.Range("D16").Value = 1.05402443947486
$ws.Range("E16").Value = 1.067962541790606
$ws.Range("F16").Value = 1.074971103989677
$ws.Range("I16").Value = 1.05165202021793
$ws.Range("J16").Value = 1.060485793352361
$ws.Range("K16").Value = 1.057398922735836
$ws.Range("L16").Value = 1.071289887693257
$ws.Range("M16").Value = 1.078275254278412
$ws.Range("N16").Value = 1.023738629025845

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.054743267398088
$ws.Range("D17").Value = 1.05436844815731
$ws.Range("E17").Value = 1.068408896726027
$ws.Range("F17").Value = 1.075446456840004
$ws.Range("I17").Value = 1.051810687345711
$ws.Range("J17").Value = 1.060801591957335
$ws.Range("K17").Value = 1.057671499547058
$ws.Range("L17").Value = 1.071665550331707
$ws.Range("M17").Value = 1.078680353486087
$ws.Range("N17").Value = 1.023849188688367

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.055006227696779
$ws.Range("D18").Value = 1.05456907332892
$ws.Range("E18").Value = 1.068669265452977
$ws.Range("F18").Value = 1.075723755330138
$ws.Range("I18").Value = 1.05190307105741
$ws.Range("J18").Value = 1.060985700490397
$ws.Range("K18").Value = 1.05783038403128
$ws.Range("L18").Value = 1.071884616456045
$ws.Range("M18").Value = 1.078916606944004
$ws.Range("N18").Value = 1.023913614411843

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.055095885916459
$ws.Range("D19").Value = 1.05463747638221
$ws.Range("E19").Value = 1.068758047458856
$ws.Range("F19").Value = 1.075818312627965
$ws.Range("I19").Value = 1.051934543693127
$ws.Range("J19").Value = 1.061048461262299
$ws.Range("K19").Value = 1.057884541755827
$ws.Range("L19").Value = 1.071959303637948
$ws.Range("M19").Value = 1.078997157530933
$ws.Range("N19").Value = 1.023935571433345

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.054694895705096
$ws.Range("D20").Value = 1.054331542313655
$ws.Range("E20").Value = 1.068361005270504
$ws.Range("F20").Value = 1.075395452552729
$ws.Range("I20").Value = 1.051793680846479
$ws.Range("J20").Value = 1.060767719222153
$ws.Range("K20").Value = 1.057642265496742
$ws.Range("L20").Value = 1.071625250624962
$ws.Range("M20").Value = 1.078636893706723
$ws.Range("N20").Value = 1.02383733308263

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.053390921853248
$ws.Range("D21").Value = 1.053336575579088
$ws.Range("E21").Value = 1.067070384831042
$ws.Range("F21").Value = 1.074021080620817
$ws.Range("I21").Value = 1.051333791582962
$ws.Range("J21").Value = 1.059853922536561
$ws.Range("K21").Value = 1.056853366571002
$ws.Range("L21").Value = 1.070538603571415
$ws.Range("M21").Value = 1.077465237976466
$ws.Range("N21").Value = 1.023517223183616

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.052570810379477
$ws.Range("D22").Value = 1.052710732346287
$ws.Range("E22").Value = 1.066259068729283
$ws.Range("F22").Value = 1.073157247496061
$ws.Range("I22").Value = 1.051043158071513
$ws.Range("J22").Value = 1.059278551287191
$ws.Range("K22").Value = 1.056356403481299
$ws.Range("L22").Value = 1.069854912934745
$ws.Range("M22").Value = 1.076728255726568
$ws.Range("N22").Value = 1.023315398584836

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.053005591411328
$ws.Range("D23").Value = 1.053042529319513
$ws.Range("E23").Value = 1.066689149133818
$ws.Range("F23").Value = 1.073615154874346
$ws.Range("I23").Value = 1.051197369155994
$ws.Range("J23").Value = 1.05958364558871
$ws.Range("K23").Value = 1.056619943429202
$ws.Range("L23").Value = 1.070217395463712
$ws.Range("M23").Value = 1.077118974330423
$ws.Range("N23").Value = 1.023422442665455

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.054716752872527
$ws.Range("D24").Value = 1.054348218566005
$ws.Range("E24").Value = 1.068382645303495
$ws.Range("F24").Value = 1.075418499091633
$ws.Range("I24").Value = 1.05180136585841
$ws.Range("J24").Value = 1.06078302513539
$ws.Range("K24").Value = 1.057655475429751
$ws.Range("L24").Value = 1.071643460486635
$ws.Range("M24").Value = 1.078656531415148
$ws.Range("N24").Value = 1.023842690311583

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.056701425795052
$ws.Range("D25").Value = 1.055862252125266
$ws.Range("E25").Value = 1.070348542189831
$ws.Range("F25").Value = 1.07751248777339
$ws.Range("I25").Value = 1.052495833633237
$ws.Range("J25").Value = 1.062171254515637
$ws.Range("K25").Value = 1.058853038107766
$ws.Range("L25").Value = 1.073296312014526
$ws.Range("M25").Value = 1.080439449404785
$ws.Range("N25").Value = 1.024327938530102
